# Apply crypto price/volume updates for Tue Nov  5 07:08:27 UTC 2024 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.769.97'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '2.430.35'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.93'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.00'
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.511'
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.168'
$ws.Range("E9").Value = '  +10.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.59'
$ws.Range("E12").Value = '  -5.99%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000176'
$ws.Range("E13").Value = '  +4.18%  '
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '68.662.16'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '2.876.25'
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").Value = '2.429.55'
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.47'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.93'
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +2.16%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.69'
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("D26").Value = '2.554.39'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("D29").Value = '0.0₃0814'
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '426.37'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '160.74'
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.08'
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '131.92'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0919'
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("E51").Value = '  +1.48%  '
